$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. The workbook carries two unused, hidden built-in cell styles left over
#    from an old hyperlink template. Drop them along with the fonts they
#    alone would otherwise keep alive.
# ---------------------------------------------------------------------------
$wb.Styles.Item("Followed Hyperlink").Delete()
$wb.Styles.Item("Hyperlink").Delete()

# ---------------------------------------------------------------------------
# 1. Remove the now-unused trailing blank rows (6:9) from the sheet.
# ---------------------------------------------------------------------------
$ws.Rows("6:9").Delete()

# ---------------------------------------------------------------------------
# 2. Row 2 currently has no explicit cell style (general formatting). Copy the
#    formats from row 1 / row 3 (bordered, wrapped Times New Roman cells) so
#    row 2 ends up visually identical to the other data rows.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("C1").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("D1").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E1").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Replace the stock-ticker table contents (3 companies -> 3 new companies).
# ---------------------------------------------------------------------------

# Row 1: NRG -> QEP
$ws.Range("A1").Value = "QEP"
$ws.Range("B1").Value = "QEP Resources Inc"
$c = $ws.Range("B1").Characters(1, 3)
$c.Font.Name = "Times New Roman"
$c.Font.Size = 11.5
$c = $ws.Range("B1").Characters(4, 14)
$c.Font.Name = "Times New Roman"
$c.Font.Size = 11.5
$ws.Range("C1").Value = 1.56
$ws.Range("D1").Value = "377.80M"
$ws.Range("E1").Value = 0.96299999999999997
$ws.Range("E1").NumberFormat = "0.000"

# Row 2: VNO -> NOG
$ws.Range("A2").Value = "NOG"
$ws.Range("B2").Value = "Northern Oil & Gas Inc"
$c = $ws.Range("B2").Characters(1, 8)
$c.Font.Name = "Times New Roman"
$c.Font.Size = 11.5
$c = $ws.Range("B2").Characters(9, 4)
$c.Font.Name = "Times New Roman"
$c.Font.Size = 11.5
$c = $ws.Range("B2").Characters(13, 2)
$c.Font.Name = "Times New Roman"
$c.Font.Size = 11.5
$c = $ws.Range("B2").Characters(15, 8)
$c.Font.Name = "Times New Roman"
$c.Font.Size = 11.5
$ws.Range("C2").Value = 1.06
$ws.Range("D2").Value = "430.14M"
$ws.Range("E2").Value = 1.2769999999999999

# Row 3: MGM -> CPE
$ws.Range("A3").Value = "CPE"
$ws.Range("B3").Value = "Callon Petroleum Co"
$ws.Range("C3").Value = 1.42
$ws.Range("D3").Value = "563.74M"
$ws.Range("E3").Value = 1.5780000000000001

# ---------------------------------------------------------------------------
# 4. Row heights: all three data rows are now a uniform 45pt.
# ---------------------------------------------------------------------------
$ws.Rows(1).RowHeight = 45
$ws.Rows(2).RowHeight = 45
$ws.Rows(3).RowHeight = 45

# ---------------------------------------------------------------------------
# 5. Column width tweaks.
# ---------------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = 9.998697916666666
$ws.Columns(2).ColumnWidth = 7.166666666666667
$ws.Columns(3).ColumnWidth = 34.830729166666664
$ws.Columns(4).ColumnWidth = 21.666666666666668
$ws.Columns(5).ColumnWidth = 10.498697916666666
$ws.Columns(6).ColumnWidth = 10.498697916666666
$ws.Columns(7).ColumnWidth = 10.998697916666666

# ---------------------------------------------------------------------------
# 6. Selection moves to C7 (mirrors the authoring tool's cursor position).
# ---------------------------------------------------------------------------
$ws.Range("C7").Select()

Write-Host "edit applied"
